$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style used by the row above (Departmental DCA row) onto the new
# "Total Forecast" row so the new cells share the same formatting.
$ws.Range("A23:F23").Style = $ws.Range("A22:F22").Style

$ws.Range("A23").Value = "Total Forecast"
$ws.Range("B23").Value = 3067
$ws.Range("C23").Value = 89
$ws.Range("D23").Value = 89809
$ws.Range("E23").Value = 87879
$ws.Range("F23").Value = 4345

$ws.Range("F25").Select()
